$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark (currently sitting after the
#    Taxila sentence) - it will be re-created later at its new location.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) "Nestorian" -> "Christian"
#    Use temporary bookmarks as barriers on each side of the word so that the
#    run containing it is not silently re-merged with its (identically
#    formatted) neighbouring runs when the text is replaced.
# ---------------------------------------------------------------------------
$searchText = $d.Content.Text
$nIdx = $searchText.IndexOf("Nestorian")
$nEnd = $nIdx + "Nestorian".Length

$d.Bookmarks.Add("ZZ_BARRIER_BEFORE", $d.Range($nIdx, $nIdx))
$d.Bookmarks.Add("ZZ_BARRIER_AFTER", $d.Range($nEnd, $nEnd))

$curText = $d.Content.Text
$nIdx2 = $curText.IndexOf("Nestorian")
$wordRange = $d.Range($nIdx2, $nIdx2 + "Nestorian".Length)
$wordRange.Text = "Christian"

$d.Bookmarks.Item("ZZ_BARRIER_BEFORE").Delete()
$d.Bookmarks.Item("ZZ_BARRIER_AFTER").Delete()

# ---------------------------------------------------------------------------
# 3) Split the "Kongphosios ..." sentence's run into two runs - one ending
#    "... inquiry and moral" and a second one starting "ity to both ends of
#    the world." - then drop the trailing double space and put the "_GoBack"
#    bookmark right after "world." (at the end of the paragraph).
# ---------------------------------------------------------------------------

# 3a) Force the run split between "moral" and "ity" using a temporary
#     bookmark placed at that boundary; the split persists even once the
#     bookmark is later removed, as long as it's still there while we do the
#     remaining edits below.
$t = $d.Content.Text
$mIdx = $t.IndexOf("morality to both ends of the world")
$splitPos = $mIdx + "moral".Length
$d.Bookmarks.Add("ZZ_SPLIT", $d.Range($splitPos, $splitPos))

# 3b) Insert the real "_GoBack" bookmark right after "...world." (before the
#     two trailing spaces).
$t2 = $d.Content.Text
$sIdx = $t2.IndexOf("ity to both ends of the world.")
$endPos = $sIdx + "ity to both ends of the world.".Length
$d.Bookmarks.Add("_GoBack", $d.Range($endPos, $endPos))

# 3c) Remove the two trailing spaces that used to follow "world." (now found
#     right after the "_GoBack" bookmark).
$t3 = $d.Content.Text
$sIdx3 = $t3.IndexOf("ity to both ends of the world.")
$endPos3 = $sIdx3 + "ity to both ends of the world.".Length
$trailingRange = $d.Range($endPos3, $endPos3 + 2)
if ($trailingRange.Text -eq "  ") {
    $trailingRange.Text = ""
}

# 3d) Remove the temporary split-marker bookmark, now that all the other
#     edits around it are done; the run split it enforced remains in place.
$d.Bookmarks.Item("ZZ_SPLIT").Delete()
